$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.345.33'
$ws.Range('E2').Value = '  +0.42%  '

$ws.Range('D3').Value = '3.330.30'
$ws.Range('E3').Value = '  -0.02%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.35'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.67%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.18'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.41%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.46%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '3.325.16'
$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.173'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +7.37%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.635'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.27%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.18'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.13%  '

$ws.Range('E13').Value = '  +2.26%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.08'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.08%  '

$ws.Range('D15').Value = '3.856.93'
$ws.Range('E15').Value = '  -0.16%  '

$ws.Range('E16').Value = '  +3.12%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.13'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.95%  '

$ws.Range('D18').Value = '3.315.63'
$ws.Range('E18').Value = '  -0.38%  '

$ws.Range('D19').Value = '64.282.45'
$ws.Range('E19').Value = '  +0.48%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.74'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.987'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.60%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '452.39'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +6.95%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.01'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +5.45%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.07'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '87.64'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.48%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.91'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +5.58%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.87'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.41%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.55'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.53%  '

$ws.Range('E29').Value = '  +0.30%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.93'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.42%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.51'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.36%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.41'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.34%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '62.14'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.91%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '572.53'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.66%  '

$ws.Range('E35').Value = '  +0.13%  '

$ws.Range('E36').Value = '  +0.07%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.142'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.67%  '

$ws.Range('E38').Value = '  -0.11%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.34'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.06%  '

$ws.Range('E40').Value = '  +0.92%  '

$ws.Range('D41').Value = '0.0₃0728'
$ws.Range('E41').Value = '  -2.38%  '

$ws.Range('D42').Value = '3.068.61'
$ws.Range('E42').Value = '  -0.66%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0414'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.49%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.74'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.32%  '

$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.46'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.64%  '

$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.19'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.41%  '

$ws.Range('E47').Value = '  +4.04%  '

$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.00%  '

$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '142.16'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.94%  '

$ws.Range('E50').Value = '  -2.15%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.14'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.25%  '
